$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-6, columns A-D
$data = @(
    @(1, 1, 5, 5),
    @(2, 1, 10, 10),
    @(5, 1, 15, 16),
    @(3, 2, 5, 5),
    @(4, 3, 5, 5)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
